# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 384, pushing the existing
# rows 384:400 down to 385:401.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 384; this shifts rows 384-400 down to 385-401
# and keeps row formatting (e.g. the date style on column D) intact.
$ws.Rows("384:384").Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(384, 1).Value  = 11
$ws.Cells.Item(384, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(384, 3).Value  = "Bíobío"
$ws.Cells.Item(384, 4).Value  = 44939
$ws.Cells.Item(384, 5).Value  = 8
$ws.Cells.Item(384, 6).Value  = "Fruta"
$ws.Cells.Item(384, 7).Value  = 100102
$ws.Cells.Item(384, 8).Value  = "Cítricos"
$ws.Cells.Item(384, 9).Value  = 100102005
$ws.Cells.Item(384, 10).Value = "Naranja"
$ws.Cells.Item(384, 11).Value = "Valencia"
$ws.Cells.Item(384, 12).Value = "Primera"
$ws.Cells.Item(384, 13).Value = 300
$ws.Cells.Item(384, 14).Value = 11000
$ws.Cells.Item(384, 15).Value = 12000
$ws.Cells.Item(384, 16).Value = 11500
$ws.Cells.Item(384, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(384, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(384, 19).Value = 767
$ws.Cells.Item(384, 20).Value = 15
